# Regenerate the "K" (strikeouts) column (column G) values for the
# armstrong_shawn save_data sheet, replacing the old Strike# derived
# values with the recalculated K values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K), rows 2 through 57, in order.
$kValues = @(3,1,2,2,1,0,3,2,1,0,2,1,4,0,1,1,3,1,2,4,2,1,0,3,0,2,1,1,3,1,0,1,1,0,3,0,0,1,2,2,2,1,0,3,0,1,1,0,2,0,1,0,1,0,0,2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
